# Update the NATMI ligand/receptor TPM-derived metrics for Cxcl12-Itgb3
# (new TPM normalization run produced different expression/specificity
# numbers). Only numeric metric columns E:T for data rows 2-17 change;
# the sending/target cluster, ligand/receptor symbol columns (A:D), and
# the header row are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 146.620486
$ws.Range("H2").Value = 439.861458
$ws.Range("I2").Value = 0.3983053592962091
$ws.Range("J2").Value = 0.3983053592962091
$ws.Range("M2").Value = 2.598166333333333
$ws.Range("N2").Value = 7.794499
$ws.Range("O2").Value = 0.3466013321552429
$ws.Range("P2").Value = 0.3466013321552429
$ws.Range("Q2").Value = 380.9444105021713
$ws.Range("R2").Value = 3428.499694519542
$ws.Range("S2").Value = 0.1380531681366387
$ws.Range("T2").Value = 0.1380531681366387
$ws.Range("G3").Value = 146.620486
$ws.Range("H3").Value = 439.861458
$ws.Range("I3").Value = 0.3983053592962091
$ws.Range("J3").Value = 0.3983053592962091
$ws.Range("M3").Value = 4.333403333333333
$ws.Range("O3").Value = 0.5780859172985858
$ws.Range("P3").Value = 0.5780859172985858
$ws.Range("Q3").Value = 635.3657027673532
$ws.Range("R3").Value = 5718.291324906179
$ws.Range("S3").Value = 0.2302547189936918
$ws.Range("T3").Value = 0.2302547189936918
$ws.Range("G4").Value = 146.620486
$ws.Range("H4").Value = 439.861458
$ws.Range("I4").Value = 0.3983053592962091
$ws.Range("J4").Value = 0.3983053592962091
$ws.Range("M4").Value = 0.4692043333333333
$ws.Range("N4").Value = 1.407613
$ws.Range("O4").Value = 0.06259293136852516
$ws.Range("P4").Value = 0.06259293136852516
$ws.Range("Q4").Value = 68.79496738663933
$ws.Range("R4").Value = 619.1547064797539
$ws.Range("S4").Value = 0.02493110001814337
$ws.Range("T4").Value = 0.02493110001814337
$ws.Range("G5").Value = 146.620486
$ws.Range("H5").Value = 439.861458
$ws.Range("I5").Value = 0.3983053592962091
$ws.Range("J5").Value = 0.3983053592962091
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.09534933333333333
$ws.Range("N5").Value = 0.286048
$ws.Range("O5").Value = 0.01271981917764605
$ws.Range("P5").Value = 0.01271981917764604
$ws.Range("Q5").Value = 13.98016559310933
$ws.Range("R5").Value = 125.821490337984
$ws.Range("S5").Value = 0.005066372147735119
$ws.Range("T5").Value = 0.005066372147735118
$ws.Range("I6").Value = 0.534552907532962
$ws.Range("J6").Value = 0.5345529075329621
$ws.Range("M6").Value = 2.598166333333333
$ws.Range("N6").Value = 7.794499
$ws.Range("O6").Value = 0.3466013321552429
$ws.Range("P6").Value = 0.3466013321552429
$ws.Range("Q6").Value = 511.2533323733866
$ws.Range("R6").Value = 4601.279991360479
$ws.Range("S6").Value = 0.185276749858383
$ws.Range("T6").Value = 0.185276749858383
$ws.Range("I7").Value = 0.534552907532962
$ws.Range("J7").Value = 0.5345529075329621
$ws.Range("M7").Value = 4.333403333333333
$ws.Range("O7").Value = 0.5780859172985858
$ws.Range("P7").Value = 0.5780859172985858
$ws.Range("Q7").Value = 852.7040267827122
$ws.Range("R7").Value = 7674.33624104441
$ws.Range("S7").Value = 0.3090175078958185
$ws.Range("T7").Value = 0.3090175078958186
$ws.Range("I8").Value = 0.534552907532962
$ws.Range("J8").Value = 0.5345529075329621
$ws.Range("M8").Value = 0.4692043333333333
$ws.Range("N8").Value = 1.407613
$ws.Range("O8").Value = 0.06259293136852516
$ws.Range("P8").Value = 0.06259293136852516
$ws.Range("Q8").Value = 92.32752957465256
$ws.Range("R8").Value = 830.947766171873
$ws.Range("S8").Value = 0.03345923345405627
$ws.Range("T8").Value = 0.03345923345405628
$ws.Range("I9").Value = 0.534552907532962
$ws.Range("J9").Value = 0.5345529075329621
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.09534933333333333
$ws.Range("N9").Value = 0.286048
$ws.Range("O9").Value = 0.01271981917764605
$ws.Range("P9").Value = 0.01271981917764604
$ws.Range("Q9").Value = 18.76233395100089
$ws.Range("R9").Value = 168.861005559008
$ws.Range("S9").Value = 0.006799416324704224
$ws.Range("T9").Value = 0.006799416324704224
$ws.Range("G10").Value = 24.174389
$ws.Range("H10").Value = 72.523167
$ws.Range("I10").Value = 0.0656715098899026
$ws.Range("J10").Value = 0.0656715098899026
$ws.Range("M10").Value = 2.598166333333333
$ws.Range("N10").Value = 7.794499
$ws.Range("O10").Value = 0.3466013321552429
$ws.Range("P10").Value = 0.3466013321552429
$ws.Range("Q10").Value = 62.80908362870367
$ws.Range("R10").Value = 565.2817526583331
$ws.Range("S10").Value = 0.02276183281248645
$ws.Range("T10").Value = 0.02276183281248645
$ws.Range("G11").Value = 24.174389
$ws.Range("H11").Value = 72.523167
$ws.Range("I11").Value = 0.0656715098899026
$ws.Range("J11").Value = 0.0656715098899026
$ws.Range("M11").Value = 4.333403333333333
$ws.Range("O11").Value = 0.5780859172985858
$ws.Range("P11").Value = 0.5780859172985858
$ws.Range("Q11").Value = 104.7573778738967
$ws.Range("R11").Value = 942.8164008650699
$ws.Range("S11").Value = 0.03796377503508749
$ws.Range("T11").Value = 0.03796377503508749
$ws.Range("G12").Value = 24.174389
$ws.Range("H12").Value = 72.523167
$ws.Range("I12").Value = 0.0656715098899026
$ws.Range("J12").Value = 0.0656715098899026
$ws.Range("M12").Value = 0.4692043333333333
$ws.Range("N12").Value = 1.407613
$ws.Range("O12").Value = 0.06259293136852516
$ws.Range("P12").Value = 0.06259293136852516
$ws.Range("Q12").Value = 11.34272807448567
$ws.Range("R12").Value = 102.084552670371
$ws.Range("S12").Value = 0.004110572311406094
$ws.Range("T12").Value = 0.004110572311406094
$ws.Range("G13").Value = 24.174389
$ws.Range("H13").Value = 72.523167
$ws.Range("I13").Value = 0.0656715098899026
$ws.Range("J13").Value = 0.0656715098899026
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.09534933333333333
$ws.Range("N13").Value = 0.286048
$ws.Range("O13").Value = 0.01271981917764605
$ws.Range("P13").Value = 0.01271981917764604
$ws.Range("Q13").Value = 2.305011874890667
$ws.Range("R13").Value = 20.745106874016
$ws.Range("S13").Value = 0.000835329730922555
$ws.Range("T13").Value = 0.0008353297309225549
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.541205
$ws.Range("H14").Value = 1.623615
$ws.Range("I14").Value = 0.001470223280926138
$ws.Range("J14").Value = 0.001470223280926138
$ws.Range("M14").Value = 2.598166333333333
$ws.Range("N14").Value = 7.794499
$ws.Range("O14").Value = 0.3466013321552429
$ws.Range("P14").Value = 0.3466013321552429
$ws.Range("Q14").Value = 1.406140610431667
$ws.Range("R14").Value = 12.655265493885
$ws.Range("S14").Value = 0.0005095813477346513
$ws.Range("T14").Value = 0.0005095813477346513
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.541205
$ws.Range("H15").Value = 1.623615
$ws.Range("I15").Value = 0.001470223280926138
$ws.Range("J15").Value = 0.001470223280926138
$ws.Range("M15").Value = 4.333403333333333
$ws.Range("O15").Value = 0.5780859172985858
$ws.Range("P15").Value = 0.5780859172985858
$ws.Range("Q15").Value = 2.345259551016666
$ws.Range("R15").Value = 21.10733595915
$ws.Range("S15").Value = 0.0008499153739879229
$ws.Range("T15").Value = 0.0008499153739879229
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.541205
$ws.Range("H16").Value = 1.623615
$ws.Range("I16").Value = 0.001470223280926138
$ws.Range("J16").Value = 0.001470223280926138
$ws.Range("M16").Value = 0.4692043333333333
$ws.Range("N16").Value = 1.407613
$ws.Range("O16").Value = 0.06259293136852516
$ws.Range("P16").Value = 0.06259293136852516
$ws.Range("Q16").Value = 0.2539357312216667
$ws.Range("R16").Value = 2.285421580995
$ws.Range("S16").Value = [double]"9.202558491941763e-05"
$ws.Range("T16").Value = [double]"9.202558491941763e-05"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.541205
$ws.Range("H17").Value = 1.623615
$ws.Range("I17").Value = 0.001470223280926138
$ws.Range("J17").Value = 0.001470223280926138
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.09534933333333333
$ws.Range("N17").Value = 0.286048
$ws.Range("O17").Value = 0.01271981917764605
$ws.Range("P17").Value = 0.01271981917764604
$ws.Range("Q17").Value = 0.05160353594666667
$ws.Range("R17").Value = 0.46443182352
$ws.Range("S17").Value = [double]"1.870097428414598e-05"
$ws.Range("T17").Value = [double]"1.870097428414598e-05"
